# Memory_map.xlsx -- "Debut du code du chauffage / Ajout du stockage des temperature"
#
# Adds a new "temperature" block to the memory map table:
#   - fills the previously-blank separator row (row 10) with
#     Begin_temperature / 0x0011- 0x0013
#   - inserts 6 new rows below it and fills the first two with
#     Actual_temperature / 0x0014 - 0x0016 and End_temperature / 0x0017 - 0x0019
#   - this pushes the "Gestion des utilisateurs" / "Door Code Pin" /
#     "Mesure de Temperature" blocks further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old E17:F18 merge (bottom spacer rows, holding the 0x10000 address) sits
# right where rows need to be inserted; unmerge it first so the insert below
# doesn't drag a merge across rows it shouldn't.
$ws.Range("E17:F18").UnMerge()

# Make room for the new temperature entries: insert 6 blank rows starting at
# row 11 (everything from the old row 11 onward shifts down to row 17+).
$ws.Rows("11:16").Insert()

# The insert copies row 10's formatting into the new rows; strip that back to
# a clean, unformatted block since only row 10 keeps the bordered/centered style.
$ws.Range("E11:F16").Clear()

# Begin_temperature entry reuses the already-styled (centered) row 10 cells.
$ws.Range("E10").Value = "Begin_temperature"
$ws.Range("G10").Value = "0x0011- 0x0013"

# Address-range column filled in first for the next two rows ...
$ws.Range("G11").Value = "0x0014 - 0x0016"
$ws.Range("G12").Value = "0x0017 - 0x0019"

# ... then the matching labels.
$ws.Range("E11").Value = "Actual_temperature"
$ws.Range("E12").Value = "End_temperature"

# Leave the selection where the author ended up before saving.
$ws.Range("I14").Select()
